# Whitelist and blacklist of classes.
# Reclassify the "Primary topic" / "Secondary topic" values for a handful
# of rows, and update the sheet's frozen-pane / selection view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Re-label Primary topic (col B) / Secondary topic (col C) for rows 58, 59, 61 ---
# Row 58: Primary "Inclusion, Diversity" -> "Leadership content"; Secondary "Smoke-free vision" -> "Inclusion, Diversity"
$ws.Range("B58").Value = "Leadership content"
$ws.Range("C58").Value = "Inclusion, Diversity"

# Row 59: Primary "Sustainability" -> "Leadership content"; Secondary "Investor Relations" -> "Sustainability"
$ws.Range("B59").Value = "Leadership content"
$ws.Range("C59").Value = "Sustainability"

# Row 61: Primary "Inclusion, Diversity" -> "Leadership content"; Secondary "Smoke-free vision" -> "Inclusion, Diversity"
$ws.Range("B61").Value = "Leadership content"
$ws.Range("C61").Value = "Inclusion, Diversity"

# --- Update the view state: frozen pane top-left cell and active selection ---
$ws.Range("A46").Select()
$excel.ActiveWindow.Panes.Item(1).FreezePanes = $false
$ws.Application.ActiveWindow.SplitColumn = 0
$ws.Application.ActiveWindow.SplitRow = 1
$ws.Application.ActiveWindow.FreezePanes = $true

$ws.Range("B61").Select()
